# Big updates to distractmaps
# Replace the "space" shared-string marker in column C with a numeric
# 0/1 flag: 1 where the cell previously read "space", 0 everywhere else
# (rows 2-31). Row 1 keeps its "corrAns" header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq "space") {
        $cell.Value = 1
    } else {
        $cell.Value = 0
    }
}

# Move the active selection, matching the author's final cursor position.
$ws.Range("F13").Select()
